# Update "Balance Sheet", update "Income Statement", and remove the
# "Cash Flow Statement" sheet entirely (persistent-upload / data
# directory restructure commit folded the cash-flow view into the
# other two statements with refreshed FY2024/FY2023 figures).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------
# Balance Sheet
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Balance Sheet")

$ws1.Range("B1").Value = "Current Year (2024)"
$ws1.Range("C1").Value = "Previous Year (2023) "

$ws1.Range("A2").Value = " Cash and Cash Equivalents"
$ws1.Range("B2").Value = "'1,000"
$ws1.Range("C2").Value = "'1,000                "

$ws1.Range("A3").Value = " Accounts Receivable"
$ws1.Range("B3").Value = "'11,987,605.97"
$ws1.Range("C3").Value = "'10,711,454.12        "

$ws1.Range("A4").Value = " Property, Plant and Equipment"
$ws1.Range("B4").Value = "'3,489,523.92"
$ws1.Range("C4").Value = "'3,494,523.92         "

$ws1.Range("A5").Value = " Total Assets"
$ws1.Range("B5").Value = "'14,355,193.96"
$ws1.Range("C5").Value = "'13,424,369.47        "

$ws1.Range("A6").Value = " Accounts Payable"
$ws1.Range("B6").Value = "'-12,443,892.15"
$ws1.Range("C6").Value = "'-10,979,515.78       "

$ws1.Range("A7").Value = " Retained Earnings"
$ws1.Range("B7").Value = "'-2,444,853.69"
$ws1.Range("C7").Value = "'-2,741,596.38        "

$ws1.Range("A8").Value = " Total Equity & Liabilities"
$ws1.Range("B8").Value = "'-14,888,745.84"
$ws1.Range("C8").Value = "'13,721,112.16        "

# ---------------------------------------------------------------
# Income Statement
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Income Statement")

$ws2.Range("B1").Value = "Current Year (2024)"
$ws2.Range("C1").Value = "Previous Year (2023) "

$ws2.Range("A2").Value = " Revenue"
$ws2.Range("B2").Value = "'-1,276,151.85"
$ws2.Range("C2").Value = "'-1,727,145.61        "

$ws2.Range("A3").Value = " Cost of Goods Sold (COGS)"
$ws2.Range("B3").Value = "'-367,148.33"
$ws2.Range("C3").Value = "'-428,513.69          "

$ws2.Range("A4").Value = " General and Administrative Expenses"
$ws2.Range("B4").Value = "'937,434.64"
$ws2.Range("C4").Value = "'1,105,786.47         "

$ws2.Range("A5").Value = " Profit Before Tax"
$ws2.Range("B5").Value = "N/A"
$ws2.Range("C5").Value = "N/A                  "

$ws2.Range("A6").Value = " Net Profit"
$ws2.Range("B6").Value = "N/A"
$ws2.Range("C6").Value = "N/A                  "

# ---------------------------------------------------------------
# Drop the Cash Flow Statement sheet (merged away in this revision)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Cash Flow Statement")
$ws3.Delete()

# Restore the original active sheet/selection
$ws1.Select()
$ws1.Range("A1").Select()
